$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Calculated PID" block in columns AD:AE (mirrors the existing
#     "Ess values" / "PD" block in W:X and "Ess values" / "PID" block in AA:AB) ---

# Header
$ws.Range("AD2").Value = "Calculated PID"

# Row 3 - sub headers
$ws.Range("AD3").Value = "Ess values"
$ws.Range("AE3").Value = "PID"

# Row 4 - Amplitude label + first sample
$ws.Range("AD4").Value = "Amplitude"
$ws.Range("AE4").Value = 1.3339799999999999

# Rows 5-7 - more amplitude samples
$ws.Range("AE5").Value = 1.2460899999999999
$ws.Range("AE6").Value = 1.2460899999999999
$ws.Range("AE7").Value = 1.3339799999999999

# Row 8 - extra sample value beneath the amplitude list
$ws.Range("AE8").Value = 0.71875

# Row 9 - extra sample value
$ws.Range("AE9").Value = 0.89453099999999997

# Row 10 - extra sample value
$ws.Range("AE10").Value = 0.54296900000000003

# Row 11 - Average label + AVERAGE formula
$ws.Range("AD11").Value = "Average"
$ws.Range("AE11").Formula = "=AVERAGE(AE4:AE10)"

# Row 12 - Setpoint label + value
$ws.Range("AD12").Value = "Setpoint = "
$ws.Range("AE12").Value = 0

# Row 13 - Overshoot label + ABS formula
$ws.Range("AD13").Value = "Overshoot"
$ws.Range("AE13").Formula = "=ABS(AE12-AE11)"

# Column AD width (matches the bestFit/customWidth column sizing used elsewhere
# on the sheet); the nearest width reproducible through this object model.
$ws.Columns("AD").ColumnWidth = 19.75

# Update the view state to match where the author left the selection/scroll
$excel.ActiveWindow.ScrollColumn = 23
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AD15").Select()
